$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H15").Value = 1622.5358
$ws.Range("I15").Value = 1622.5358
$ws.Range("J15").Value = 0
$ws.Range("K15").Value = 4867.607400000001
$ws.Range("L15").Value = 0
$ws.Range("M15").Value = -4698.607400000001

$ws.Range("H62").Value = 4064.7
$ws.Range("I62").Value = 3664.4285
$ws.Range("J62").Value = 4998.6665
$ws.Range("K62").Value = 3664.4285
$ws.Range("L62").Value = 4998.6665
$ws.Range("M62").Value = -3040.4285
$ws.Range("N62").Value = -6246.6665

$ws.Range("H65").Value = 4064.7
$ws.Range("I65").Value = 3664.4285
$ws.Range("J65").Value = 4998.6665
$ws.Range("K65").Value = 18322.1425
$ws.Range("L65").Value = 24993.3325
$ws.Range("M65").Value = -15202.1425
$ws.Range("N65").Value = -31233.3325

$ws.Range("H76").Value = 10000
$ws.Range("I76").Value = 0
$ws.Range("J76").Value = 10000
$ws.Range("K76").Value = 0
$ws.Range("L76").Value = 10000
$ws.Range("N76").Value = -10630
$ws.Range("M76").ClearContents()

$ws.Range("H79").Value = 10000
$ws.Range("I79").Value = 0
$ws.Range("J79").Value = 10000
$ws.Range("K79").Value = 0
$ws.Range("L79").Value = 10000
$ws.Range("N79").Value = -12184
$ws.Range("M79").ClearContents()

$ws.Range("H87").Value = 21363.637
$ws.Range("I87").Value = 0
$ws.Range("J87").Value = 21363.637
$ws.Range("K87").Value = 0
$ws.Range("L87").Value = 21363.637
$ws.Range("N87").Value = -23859.637

$ws.Range("H90").Value = 21363.637
$ws.Range("I90").Value = 0
$ws.Range("J90").Value = 21363.637
$ws.Range("K90").Value = 0
$ws.Range("L90").Value = 64090.91099999999
$ws.Range("N90").Value = -76570.91099999999

$ws.Range("H106").Value = 1462.875
$ws.Range("I106").Value = 1314.7142
$ws.Range("J106").Value = 2500
$ws.Range("K106").Value = 1314.7142
$ws.Range("L106").Value = 2500
$ws.Range("M106").Value = -683.7141999999999
$ws.Range("N106").Value = -3762

$ws.Range("H137").Value = 3209.7046
$ws.Range("I137").Value = 3033.7368
$ws.Range("J137").Value = 4324.1665
$ws.Range("K137").Value = 9101.2104
$ws.Range("L137").Value = 12972.4995
$ws.Range("M137").Value = -6551.2104
$ws.Range("N137").Value = -18072.4995

$ws.Range("H141").Value = 4420.757
$ws.Range("I141").Value = 2143.25
$ws.Range("J141").Value = 18996.8
$ws.Range("K141").Value = 6429.75
$ws.Range("L141").Value = 56990.39999999999
$ws.Range("M141").Value = -1249.75
$ws.Range("N141").Value = -67350.39999999999

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 25933.615
$ws.Range("I2").Value = 54457.332
$ws.Range("J2").Value = 1484.7142
$ws.Range("K2").Value = 54457.332
$ws.Range("L2").Value = 1484.7142
$ws.Range("M2").Value = -54344.332
$ws.Range("N2").Value = -1710.7142

$ws.Range("H32").Value = 6822.759
$ws.Range("I32").Value = 7314.08
$ws.Range("J32").Value = 3752
$ws.Range("K32").Value = 7314.08
$ws.Range("L32").Value = 3752
$ws.Range("M32").Value = -7027.08
$ws.Range("N32").Value = -4326

$ws.Range("H74").Value = 1474
$ws.Range("I74").Value = 1237.7273
$ws.Range("J74").Value = 2123.75
$ws.Range("K74").Value = 1237.7273
$ws.Range("L74").Value = 2123.75
$ws.Range("M74").Value = -363.7273
$ws.Range("N74").Value = -3871.75

$ws.Range("H77").Value = 1474
$ws.Range("I77").Value = 1237.7273
$ws.Range("J77").Value = 2123.75
$ws.Range("K77").Value = 6188.636500000001
$ws.Range("L77").Value = 10618.75
$ws.Range("M77").Value = -1820.636500000001
$ws.Range("N77").Value = -19354.75

$ws.Range("H97").Value = 18522872
$ws.Range("I97").Value = 1529
$ws.Range("J97").Value = 55565556
$ws.Range("K97").Value = 1529
$ws.Range("L97").Value = 55565556
$ws.Range("M97").Value = -1033
$ws.Range("N97").Value = -55566548

$ws.Range("H102").Value = 18520464
$ws.Range("I102").Value = 2028.0714
$ws.Range("J102").Value = 83334990
$ws.Range("K102").Value = 2028.0714
$ws.Range("L102").Value = 83334990
$ws.Range("M102").Value = -406.0714
$ws.Range("N102").Value = -83338234

$ws.Range("H116").Value = 25933.615
$ws.Range("I116").Value = 54457.332
$ws.Range("J116").Value = 1484.7142
$ws.Range("K116").Value = 54457.332
$ws.Range("L116").Value = 1484.7142
$ws.Range("M116").Value = -52163.332
$ws.Range("N116").Value = -6072.7142

$ws.Range("H132").Value = 1899.0344
$ws.Range("I132").Value = 1942.2727
$ws.Range("J132").Value = 1763.1428
$ws.Range("K132").Value = 5826.8181
$ws.Range("L132").Value = 5289.428400000001
$ws.Range("M132").Value = -3296.8181
$ws.Range("N132").Value = -10349.4284

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 25933.615
$ws.Range("I3").Value = 54457.332
$ws.Range("J3").Value = 1484.7142
$ws.Range("K3").Value = 54457.332
$ws.Range("L3").Value = 1484.7142
$ws.Range("M3").Value = -54343.332
$ws.Range("N3").Value = -1712.7142

$ws.Range("H11").Value = 1562.75
$ws.Range("I11").Value = 127.75
$ws.Range("J11").Value = 2997.75
$ws.Range("K11").Value = 127.75
$ws.Range("L11").Value = 2997.75
$ws.Range("M11").Value = 12.25
$ws.Range("N11").Value = -3277.75

$ws.Range("H20").Value = 1610.75
$ws.Range("I20").Value = 1373.8
$ws.Range("J20").Value = 2005.6666
$ws.Range("K20").Value = 1373.8
$ws.Range("L20").Value = 2005.6666
$ws.Range("M20").Value = -1126.8
$ws.Range("N20").Value = -2499.6666

$ws.Range("H107").Value = 2039.1428
$ws.Range("I107").Value = 2059.111
$ws.Range("J107").Value = 1500
$ws.Range("K107").Value = 2059.111
$ws.Range("L107").Value = 1500
$ws.Range("M107").Value = -139.1109999999999
$ws.Range("N107").Value = -5340

$ws.Range("H134").Value = 2833.4856
$ws.Range("I134").Value = 2663.2964
$ws.Range("J134").Value = 3407.875
$ws.Range("K134").Value = 7989.889200000001
$ws.Range("L134").Value = 10223.625
$ws.Range("M134").Value = -5454.889200000001
$ws.Range("N134").Value = -15293.625

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 441.2857
$ws.Range("I22").Value = 345
$ws.Range("J22").Value = 682
$ws.Range("K22").Value = 345
$ws.Range("L22").Value = 682
$ws.Range("M22").Value = 5
$ws.Range("N22").Value = -1382

$ws.Range("H31").Value = 2072.9792
$ws.Range("I31").Value = 1833.5555
$ws.Range("J31").Value = 2791.25
$ws.Range("K31").Value = 1833.5555
$ws.Range("L31").Value = 2791.25
$ws.Range("M31").Value = -1538.5555
$ws.Range("N31").Value = -3381.25

$ws.Range("H34").Value = 2072.9792
$ws.Range("I34").Value = 1833.5555
$ws.Range("J34").Value = 2791.25
$ws.Range("K34").Value = 1833.5555
$ws.Range("L34").Value = 2791.25
$ws.Range("M34").Value = -1631.5555
$ws.Range("N34").Value = -3128.8333

$ws.Range("H86").Value = 53552.555
$ws.Range("I86").Value = 85001.5
$ws.Range("J86").Value = 28393.4
$ws.Range("K86").Value = 85001.5
$ws.Range("L86").Value = 28393.4
$ws.Range("M86").Value = -83878.5
$ws.Range("N86").Value = -30639.4

$ws.Range("H89").Value = 53552.555
$ws.Range("I89").Value = 85001.5
$ws.Range("J89").Value = 28393.4
$ws.Range("K89").Value = 425007.5
$ws.Range("L89").Value = 141967
$ws.Range("M89").Value = -419391.5
$ws.Range("N89").Value = -153199

$ws.Range("H107").Value = 9724.478999999999
$ws.Range("I107").Value = 816.86664
$ws.Range("J107").Value = 26426.25
$ws.Range("K107").Value = 816.86664
$ws.Range("L107").Value = 26426.25
$ws.Range("M107").Value = 1103.13336
$ws.Range("N107").Value = -30266.25

$ws.Range("H132").Value = 2438.1943
$ws.Range("I132").Value = 2321.7742
$ws.Range("J132").Value = 3160
$ws.Range("K132").Value = 6965.3226
$ws.Range("L132").Value = 9480
$ws.Range("M132").Value = -4435.3226
$ws.Range("N132").Value = -14540

$ws.Range("H133").Value = 80000
$ws.Range("I133").Value = 0
$ws.Range("J133").Value = 80000
$ws.Range("K133").Value = 0
$ws.Range("L133").Value = 80000
$ws.Range("N133").Value = -85060

$ws.Range("H134").Value = 5134.45
$ws.Range("I134").Value = 4609.2
$ws.Range("J134").Value = 5659.7
$ws.Range("K134").Value = 13827.6
$ws.Range("L134").Value = 16979.1
$ws.Range("M134").Value = -11292.6
$ws.Range("N134").Value = -22049.1

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H41").Value = 592.2308
$ws.Range("I41").Value = 1799.6666
$ws.Range("J41").Value = 230
$ws.Range("K41").Value = 5398.9998
$ws.Range("L41").Value = 690
$ws.Range("M41").Value = -5060.9998
$ws.Range("N41").Value = -1366

$ws.Range("H118").Value = 4327
$ws.Range("I118").Value = 981.3333
$ws.Range("J118").Value = 5999.8335
$ws.Range("K118").Value = 2943.9999
$ws.Range("L118").Value = 17999.5005
$ws.Range("M118").Value = -1700.9999
$ws.Range("N118").Value = -20485.5005

$ws.Range("H129").Value = 2109.7693
$ws.Range("I129").Value = 900
$ws.Range("J129").Value = 2210.5833
$ws.Range("K129").Value = 2700
$ws.Range("L129").Value = 6631.749899999999
$ws.Range("M129").Value = 2300
$ws.Range("N129").Value = -16631.7499

$ws.Range("H131").Value = 3253.182
$ws.Range("I131").Value = 1784.2222
$ws.Range("J131").Value = 4270.154
$ws.Range("K131").Value = 5352.6666
$ws.Range("L131").Value = 12810.462
$ws.Range("M131").Value = -312.6665999999996
$ws.Range("N131").Value = -22890.462

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H107").Value = 863.5714
$ws.Range("I107").Value = 1272
$ws.Range("J107").Value = 319
$ws.Range("K107").Value = 1272
$ws.Range("L107").Value = 319
$ws.Range("M107").Value = 648
$ws.Range("N107").Value = -4159

$ws.Range("H113").Value = 4334.2856
$ws.Range("I113").Value = 1171.25
$ws.Range("J113").Value = 8551.666999999999
$ws.Range("K113").Value = 1171.25
$ws.Range("L113").Value = 8551.666999999999
$ws.Range("M113").Value = 998.75
$ws.Range("N113").Value = -12891.667

$ws.Range("H122").Value = 3838.2144
$ws.Range("I122").Value = 2042.8462
$ws.Range("J122").Value = 5394.2
$ws.Range("K122").Value = 6128.5386
$ws.Range("L122").Value = 16182.6
$ws.Range("M122").Value = -3678.5386
$ws.Range("N122").Value = -21082.6

$ws.Range("H132").Value = 2631.0293
$ws.Range("I132").Value = 2772.2273
$ws.Range("J132").Value = 2372.1667
$ws.Range("K132").Value = 8316.6819
$ws.Range("L132").Value = 7116.500100000001
$ws.Range("M132").Value = -5786.6819
$ws.Range("N132").Value = -12176.5001

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 6440.8
$ws.Range("I7").Value = 8267.666999999999
$ws.Range("J7").Value = 5657.857
$ws.Range("K7").Value = 8267.666999999999
$ws.Range("L7").Value = 5657.857
$ws.Range("M7").Value = -8155.666999999999
$ws.Range("N7").Value = -5881.857

$ws.Range("H68").Value = 5000
$ws.Range("I68").Value = 5000
$ws.Range("J68").Value = 0
$ws.Range("K68").Value = 5000
$ws.Range("L68").Value = 0
$ws.Range("M68").Value = -4251

$ws.Range("H71").Value = 5000
$ws.Range("I71").Value = 5000
$ws.Range("J71").Value = 0
$ws.Range("K71").Value = 25000
$ws.Range("L71").Value = 0
$ws.Range("M71").Value = -21256

$ws.Range("H100").Value = 78886.60000000001
$ws.Range("I100").Value = 161615
$ws.Range("J100").Value = 6499.25
$ws.Range("K100").Value = 161615
$ws.Range("L100").Value = 6499.25
$ws.Range("M100").Value = -161074
$ws.Range("N100").Value = -7581.25

$ws.Range("H126").Value = 6440.8
$ws.Range("I126").Value = 8267.666999999999
$ws.Range("J126").Value = 5657.857
$ws.Range("K126").Value = 24803.001
$ws.Range("L126").Value = 16973.571
$ws.Range("M126").Value = -22333.001
$ws.Range("N126").Value = -21913.571

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H16").Value = 0
$ws.Range("I16").Value = 0
$ws.Range("J16").Value = 0
$ws.Range("K16").Value = 0
$ws.Range("L16").ClearContents()
$ws.Range("N16").ClearContents()

$ws.Range("H81").Value = 4834492.5
$ws.Range("I81").Value = 3515.7144
$ws.Range("J81").Value = 12349345
$ws.Range("K81").Value = 7031.4288
$ws.Range("L81").Value = 24698690
$ws.Range("M81").Value = -5970.4288
$ws.Range("N81").Value = -24700812

$ws.Range("H84").Value = 4834492.5
$ws.Range("I84").Value = 3515.7144
$ws.Range("J84").Value = 12349345
$ws.Range("K84").Value = 35157.144
$ws.Range("L84").Value = 123493450
$ws.Range("M84").Value = -29853.144
$ws.Range("N84").Value = -123504058

$ws.Range("H96").Value = 58322.11
$ws.Range("I96").Value = 86050
$ws.Range("J96").Value = 2866.3333
$ws.Range("K96").Value = 86050
$ws.Range("L96").Value = 2866.3333
$ws.Range("M96").Value = -86677
$ws.Range("N96").Value = -5612.3333
